# MahanaAamdanDeposits1YearRegularCustomer.xlsx - add T.C (Azure) / T.C (Desc.) / Error
# columns with a bold 14pt header row and a boxed-table border (medium outer edge,
# thin inner gridlines, medium line under the header and under the data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header values (columns C, D, E) ------------------------------------
$ws.Range("C1").Value = "T.C (Azure)"
$ws.Range("D1").Value = "T.C (Desc.)"
$ws.Range("E1").Value = "Error"

# ---- New data row values ------------------------------------------------------
$ws.Range("C2").Value = 119492
# D2 / E2 stay empty (styled only)

# ---- Header row formatting: bold, 14pt font ----------------------------------
$header = $ws.Range("A1:E1")
$header.Font.Bold = $true
$header.Font.Size = 14

# ---- Row heights ---------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 15.75

# ---- Table borders: medium outside border, thin inside gridlines --------------
$table = $ws.Range("A1:E2")

# Start every edge (outside + inside) as a thin continuous line ...
$table.Borders.LineStyle = 1

# ... then thicken the top/bottom/left/right OUTSIDE edges to medium.
$table.Borders.Item(8).Weight = -4138
$table.Borders.Item(9).Weight = -4138
$table.Borders.Item(7).Weight = -4138
$table.Borders.Item(10).Weight = -4138

# ---- Data row alignment: the T.C (Azure) figure is left-aligned ---------------
$ws.Range("C2").HorizontalAlignment = -4131

# ---- Selection matches the authored state (active cell C2) --------------------
$ws.Range("C2").Select()
